# Insert a new row at position 193 (pushes existing rows 193:262 down to 194:263)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(193).Insert()

# Populate the newly inserted row 193 with the new record's data
$ws.Range("A193").Value = 10
$ws.Range("B193").Value = "Vega Modelo de Temuco"
$ws.Range("C193").Value = "La Araucanía"
$ws.Range("D193").Value = 44588
$ws.Range("E193").Value = 9
$ws.Range("F193").Value = 100112044
$ws.Range("G193").Value = "Perejil"
$ws.Range("H193").Value = "Sin especificar"
$ws.Range("I193").Value = "Primera"
$ws.Range("J193").Value = 55
$ws.Range("K193").Value = 6000
$ws.Range("L193").Value = 6000
$ws.Range("M193").Value = 6000
$ws.Range("N193").Value = "$/docena de atados (3 kilos)"
$ws.Range("O193").Value = "Provincia de Cautín"
$ws.Range("P193").Value = 2000
$ws.Range("Q193").Value = 3
$ws.Range("R193").Value = "Hortaliza"
